# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# updates across the ALC, ARM, CRP, GSM, and LTW sheets (per commit diff).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3012.4583
$ws.Range("J64").Value = 3004.9473
$ws.Range("L64").Value = 3004.9473
$ws.Range("N64").Value = -3500.9473
$ws.Range("H67").Value = 3012.4583
$ws.Range("J67").Value = 3004.9473
$ws.Range("L67").Value = 3004.9473
$ws.Range("N67").Value = -4720.9473
$ws.Range("H76").Value = 92027.03999999999
$ws.Range("I76").Value = 113883.6
$ws.Range("J76").Value = 4600.8
$ws.Range("K76").Value = 113883.6
$ws.Range("L76").Value = 4600.8
$ws.Range("M76").Value = -113568.6
$ws.Range("N76").Value = -5230.8
$ws.Range("H79").Value = 92027.03999999999
$ws.Range("I79").Value = 113883.6
$ws.Range("J79").Value = 4600.8
$ws.Range("K79").Value = 113883.6
$ws.Range("L79").Value = 4600.8
$ws.Range("M79").Value = -112791.6
$ws.Range("N79").Value = -6784.8
$ws.Range("H132").Value = 3273.577
$ws.Range("I132").Value = 3077.5232
$ws.Range("K132").Value = 9232.569600000001
$ws.Range("M132").Value = -6702.569600000001
$ws.Range("H137").Value = 56847.668
$ws.Range("I137").Value = 925.6
$ws.Range("K137").Value = 2776.8
$ws.Range("M137").Value = -226.8000000000002

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26596.88
$ws.Range("I32").Value = 15044.638
$ws.Range("J32").Value = 42549.977
$ws.Range("K32").Value = 15044.638
$ws.Range("L32").Value = 42549.977
$ws.Range("M32").Value = -14757.638
$ws.Range("N32").Value = -43123.977
$ws.Range("H63").Value = 2209.9285
$ws.Range("I63").Value = 2226.077
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 2226.077
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -1540.077
$ws.Range("N63").Value = -3372
$ws.Range("H66").Value = 2209.9285
$ws.Range("I66").Value = 2226.077
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 11130.385
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -7698.385000000002
$ws.Range("N66").Value = -16864
$ws.Range("H88").Value = 68171.75
$ws.Range("I88").Value = 1757
$ws.Range("J88").Value = 74209.45
$ws.Range("K88").Value = 1757
$ws.Range("L88").Value = 74209.45
$ws.Range("M88").Value = -1351
$ws.Range("N88").Value = -75021.45
$ws.Range("H91").Value = 68171.75
$ws.Range("I91").Value = 1757
$ws.Range("J91").Value = 74209.45
$ws.Range("K91").Value = 1757
$ws.Range("L91").Value = 74209.45
$ws.Range("M91").Value = -353
$ws.Range("N91").Value = -77017.45

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2207.158
$ws.Range("I31").Value = 1237.3715
$ws.Range("J31").Value = 3750
$ws.Range("K31").Value = 1237.3715
$ws.Range("L31").Value = 3750
$ws.Range("M31").Value = -942.3715
$ws.Range("N31").Value = -4340
$ws.Range("H34").Value = 2207.158
$ws.Range("I34").Value = 1237.3715
$ws.Range("J34").Value = 3750
$ws.Range("K34").Value = 1237.3715
$ws.Range("L34").Value = 3750
$ws.Range("M34").Value = -1035.3715
$ws.Range("N34").Value = -4154
$ws.Range("H62").Value = 3522.2222
$ws.Range("I62").Value = 2950
$ws.Range("J62").Value = 3685.7144
$ws.Range("K62").Value = 2950
$ws.Range("L62").Value = 3685.7144
$ws.Range("M62").Value = -2326
$ws.Range("N62").Value = -4933.7144
$ws.Range("H65").Value = 3522.2222
$ws.Range("I65").Value = 2950
$ws.Range("J65").Value = 3685.7144
$ws.Range("K65").Value = 14750
$ws.Range("L65").Value = 18428.572
$ws.Range("M65").Value = -11630
$ws.Range("N65").Value = -24668.572
$ws.Range("H132").Value = 4511.375
$ws.Range("I132").Value = 4059.6667
$ws.Range("J132").Value = 4782.4
$ws.Range("K132").Value = 12179.0001
$ws.Range("L132").Value = 14347.2
$ws.Range("M132").Value = -9649.000100000001
$ws.Range("N132").Value = -19407.2
$ws.Range("H134").Value = 3093.0312
$ws.Range("I134").Value = 3209.2334
$ws.Range("J134").Value = 1350
$ws.Range("K134").Value = 9627.700199999999
$ws.Range("L134").Value = 4050
$ws.Range("M134").Value = -7092.700199999999
$ws.Range("N134").Value = -9120

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 33400.5
$ws.Range("J34").Value = 33400.5
$ws.Range("L34").Value = 33400.5
$ws.Range("N34").Value = -33936.5
$ws.Range("H62").Value = 30000
$ws.Range("I62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("K62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("M62").Value = -29314
$ws.Range("N62").Value = -31372
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("H65").Value = 30000
$ws.Range("I65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("K65").Value = 90000
$ws.Range("L65").Value = 90000
$ws.Range("M65").Value = -86568
$ws.Range("N65").Value = -96864
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("H68").Value = 38000
$ws.Range("J68").Value = 38000
$ws.Range("L68").Value = 38000
$ws.Range("N68").Value = -39622
$ws.Range("H69").Value = 25000
$ws.Range("J69").Value = 25000
$ws.Range("L69").Value = 25000
$ws.Range("N69").Value = -26498
$ws.Range("H70").Value = 165885360
$ws.Range("I70").Value = 414705900
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 414705900
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -414705630
$ws.Range("N70").Value = -5540
$ws.Range("H71").Value = 38000
$ws.Range("J71").Value = 38000
$ws.Range("L71").Value = 114000
$ws.Range("N71").Value = -122112
$ws.Range("H72").Value = 25000
$ws.Range("J72").Value = 25000
$ws.Range("L72").Value = 75000
$ws.Range("N72").Value = -82488
$ws.Range("H73").Value = 165885360
$ws.Range("I73").Value = 414705900
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 414705900
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -414704964
$ws.Range("N73").Value = -6872
$ws.Range("H76").Value = 33400.5
$ws.Range("J76").Value = 33400.5
$ws.Range("L76").Value = 33400.5
$ws.Range("N76").Value = -34030.5
$ws.Range("H79").Value = 33400.5
$ws.Range("J79").Value = 33400.5
$ws.Range("L79").Value = 33400.5
$ws.Range("N79").Value = -35584.5
$ws.Range("H80").Value = 5100.7144
$ws.Range("I80").Value = 8701.666999999999
$ws.Range("J80").Value = 2400
$ws.Range("K80").Value = 8701.666999999999
$ws.Range("L80").Value = 2400
$ws.Range("M80").Value = -7703.666999999999
$ws.Range("N80").Value = -4396
$ws.Range("H83").Value = 5100.7144
$ws.Range("I83").Value = 8701.666999999999
$ws.Range("J83").Value = 2400
$ws.Range("K83").Value = 43508.335
$ws.Range("L83").Value = 12000
$ws.Range("M83").Value = -38516.335
$ws.Range("N83").Value = -21984
$ws.Range("H132").Value = 3838
$ws.Range("I132").Value = 4051.7144
$ws.Range("J132").Value = 3196.8572
$ws.Range("K132").Value = 12155.1432
$ws.Range("L132").Value = 9590.571599999999
$ws.Range("M132").Value = -9625.143199999999
$ws.Range("N132").Value = -14650.5716
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1509.091
$ws.Range("J46").Value = 700
$ws.Range("L46").Value = 700
$ws.Range("N46").Value = -1076
$ws.Range("H68").Value = 387023
$ws.Range("I68").Value = 626199.9
$ws.Range("J68").Value = 4340
$ws.Range("K68").Value = 626199.9
$ws.Range("L68").Value = 4340
$ws.Range("M68").Value = -625450.9
$ws.Range("N68").Value = -5838
$ws.Range("H71").Value = 387023
$ws.Range("I71").Value = 626199.9
$ws.Range("J71").Value = 4340
$ws.Range("K71").Value = 3130999.5
$ws.Range("L71").Value = 21700
$ws.Range("M71").Value = -3127255.5
$ws.Range("N71").Value = -29188
$ws.Range("H74").Value = 11500
$ws.Range("J74").Value = 11500
$ws.Range("L74").Value = 11500
$ws.Range("N74").Value = -13496
$ws.Range("H76").Value = 33333
$ws.Range("J76").Value = 33333
$ws.Range("L76").Value = 33333
$ws.Range("N76").Value = -34009
$ws.Range("H77").Value = 11500
$ws.Range("J77").Value = 11500
$ws.Range("L77").Value = 34500
$ws.Range("N77").Value = -44484
$ws.Range("H79").Value = 33333
$ws.Range("J79").Value = 33333
$ws.Range("L79").Value = 33333
$ws.Range("N79").Value = -35673
